$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "Thomas Debiasi"
$ws.Range("B66").Value = "Thomas Debiasi | Mai una gioia"
$ws.Range("C66").Value = "Randy Cobbinah | Mai una gioia"
$ws.Range("D66").Value = "Stefano Mattioli | SdrumALA"
$ws.Range("E66").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("F66").Value = "Federico Nicolodi | U.SGUARNA"
